$d = $word.ActiveDocument

# Move to the end of the document's main story and add a new paragraph
# after the current last one ("When state passed then it is passed as
# props."). Word inherits the preceding paragraph's formatting (the
# ListParagraph style plus the existing numbered-list numPr/numId) for
# a paragraph inserted this way, so the new bullet lands in the same
# list without minting a new numbering definition.
$end = $d.Content
$end.Collapse(0)  # wdCollapseEnd
$end.InsertParagraphAfter()

$newPara = $d.Paragraphs.Last
$newPara.Range.Text = "For random id npm i uuid"
